# Add a new standard row (UWG-2) to the WiscSIMS run standards sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20: StdName, d18O, REGEX, StdType (no d13C value for this standard)
$ws.Range("A20").Value = "UWG-2"
$ws.Range("C20").Value = 5.8
$ws.Range("D20").Value = "UWG\D*2"
$ws.Range("E20").Value = "Run"

# Match the author's final cell selection left in the saved workbook
$ws.Range("F17").Select()
